$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# "Boolean" sheet: split the two combined trans CSV-list entries
# (trans/BVTQaZ/BVTQaZ.csv and trans/VTQaZ/VTQaZ.csv) into six per-mode files
# each, keeping everything else in the list in the same relative order.
# ---------------------------------------------------------------------------
$wsBool = $wb.Worksheets.Item("Boolean")

# Row 17 currently holds "trans/BVTQaZ/BVTQaZ.csv" -> replace + insert 5 more
$wsBool.Range("A17").Value2 = "trans/BVTQaZ/BVTQaZ-LDVs.csv"
$wsBool.Rows.Item(18).Insert()
$wsBool.Rows.Item(18).Insert()
$wsBool.Rows.Item(18).Insert()
$wsBool.Rows.Item(18).Insert()
$wsBool.Rows.Item(18).Insert()
$wsBool.Range("A18").Value2 = "trans/BVTQaZ/BVTQaZ-HDVs.csv"
$wsBool.Range("A19").Value2 = "trans/BVTQaZ/BVTQaZ-aircraft.csv"
$wsBool.Range("A20").Value2 = "trans/BVTQaZ/BVTQaZ-rail.csv"
$wsBool.Range("A21").Value2 = "trans/BVTQaZ/BVTQaZ-ships.csv"
$wsBool.Range("A22").Value2 = "trans/BVTQaZ/BVTQaZ-motorbikes.csv"

# After the above, "trans/VTQaZ/VTQaZ.csv" (previously row 21) is now row 26.
$wsBool.Range("A26").Value2 = "trans/VTQaZ/VTQaZ-LDVs.csv"
$wsBool.Rows.Item(27).Insert()
$wsBool.Rows.Item(27).Insert()
$wsBool.Rows.Item(27).Insert()
$wsBool.Rows.Item(27).Insert()
$wsBool.Rows.Item(27).Insert()
$wsBool.Range("A27").Value2 = "trans/VTQaZ/VTQaZ-HDVs.csv"
$wsBool.Range("A28").Value2 = "trans/VTQaZ/VTQaZ-aircraft.csv"
$wsBool.Range("A29").Value2 = "trans/VTQaZ/VTQaZ-rail.csv"
$wsBool.Range("A30").Value2 = "trans/VTQaZ/VTQaZ-ships.csv"
$wsBool.Range("A31").Value2 = "trans/VTQaZ/VTQaZ-motorbikes.csv"

# Six blank (but formatted) trailer rows after the list, rows 33-38.
$wsBool.Rows.Item(33).Font.Name = "Calibri"
$wsBool.Rows.Item(34).Font.Name = "Calibri"
$wsBool.Rows.Item(35).Font.Name = "Calibri"
$wsBool.Rows.Item(36).Font.Name = "Calibri"
$wsBool.Rows.Item(37).Font.Name = "Calibri"
$wsBool.Rows.Item(38).Font.Name = "Calibri"

# View state: scrolled down a bit, with A32 selected.
$wsBool.Activate()
$excel.ActiveWindow.ScrollRow = 10
$wsBool.Range("A32").Select()

# ---------------------------------------------------------------------------
# View / active-tab bookkeeping: "Integer" sheet keeps a remembered selection
# at A13 but is no longer the active tab; "About" becomes the active tab.
# ---------------------------------------------------------------------------
$wsInt = $wb.Worksheets.Item("Integer")
$wsInt.Range("A13").Select()

$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Activate()
